$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B cells need style s="4" (matching column A in the same row / the rest of the block);
# copy format from an existing s="4" cell before setting the value.
$ws.Range("A70").Copy()
$ws.Range("B70").PasteSpecial(-4122)
$ws.Range("B72").PasteSpecial(-4122)
$ws.Range("B73").PasteSpecial(-4122)
$ws.Range("B74").PasteSpecial(-4122)
$ws.Range("B75").PasteSpecial(-4122)
$ws.Range("B76").PasteSpecial(-4122)
$ws.Range("B77").PasteSpecial(-4122)
$ws.Range("B78").PasteSpecial(-4122)
$ws.Range("B79").PasteSpecial(-4122)
$ws.Range("B80").PasteSpecial(-4122)
$ws.Range("B81").PasteSpecial(-4122)
$ws.Range("B82").PasteSpecial(-4122)
$ws.Range("B83").PasteSpecial(-4122)
$ws.Range("B84").PasteSpecial(-4122)

# Now populate the Center of Mass values for rows 69-94 (new plot data).
$ws.Range("B69").Value = 352.239669421487
$ws.Range("B70").Value = 504.270531400966
$ws.Range("B71").Value = 565.137346938775
$ws.Range("B72").Value = 332.905671296296
$ws.Range("B73").Value = 649.214078859434
$ws.Range("B74").Value = 615.777233236801
$ws.Range("B75").Value = 532.90478079746
$ws.Range("B76").Value = 388.81663516068
$ws.Range("B77").Value = 567.150590331922
$ws.Range("B78").Value = 545.177346938775
$ws.Range("B79").Value = 421.778197857592
$ws.Range("B80").Value = 442.411672978391
$ws.Range("B81").Value = 502.466326530612
$ws.Range("B82").Value = 628.646944444444
$ws.Range("B83").Value = 565.516219723183
$ws.Range("B84").Value = 396.082647462277
$ws.Range("B85").Value = 477.320501730103
$ws.Range("B86").Value = 566.722130177514
$ws.Range("B87").Value = 257.26275510204
$ws.Range("B88").Value = 510.385637770283
$ws.Range("B89").Value = 570.082040816326
$ws.Range("B90").Value = 616.135290616941
$ws.Range("B91").Value = 320.80574845679
$ws.Range("B92").Value = 522.467959183673
$ws.Range("B93").Value = 486.636734693877
$ws.Range("B94").Value = 557.247551020408
